$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet (book-level display name, also reflected in tab)
$ws.Name = "CubeA"

# Tiny floating point corrections to existing cells (1-ULP recalculation drift)
$ws.Range("I13").Value = 0.9947087580526029
$ws.Range("O15").Value = 0.9917731065300173
$ws.Range("P15").Value = 0.9919317152281452

# New row 16 of averaged-intensity data (Gaussian Quadrature scheme export)
$ws.Range("A16").Value = 14
# Reuse the bold/bordered/centered "index column" style already applied to
# A2:A15 instead of re-declaring fonts/borders (keeps styles.xml untouched).
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9281321562220225
$ws.Range("D16").Value = 1.096363831534633
$ws.Range("E16").Value = 0.9768042872944713
$ws.Range("F16").Value = 1.014241451521715
$ws.Range("G16").Value = 0.9281321562220225
$ws.Range("H16").Value = 1.096363831534633
$ws.Range("I16").Value = 0.9635793264604577
$ws.Range("J16").Value = 1.014241451521715
$ws.Range("K16").Value = 0.9698176992801386
$ws.Range("L16").Value = 1.057161961425036
$ws.Range("M16").Value = 0.9281321562220225
$ws.Range("N16").Value = 1.036584059414552
$ws.Range("O16").Value = 1.00388543164321
$ws.Range("P16").Value = 1.002542770657524
